{"js": "// Load all paragraphs in the document body so we can locate the\n// empty list-item paragraph that sits right before the \"Push\" heading\n// (it currently has no runs, only the list/lang paragraph properties).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet target = null;\nfor (let i = 0; i < items.length - 1; i++) {\n  if (items[i].text === \"\" && items[i + 1].text === \"Push\") {\n    target = items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the empty paragraph before 'Push'.\");\n}\n\n// Insert the two runs (\"Git\" wrapped in spell-check markers, then the\n// rest of the sentence) via OOXML so the run/proofErr structure matches\n// exactly what Word itself would author.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:rPr><w:lang w:val=\"en-US\"/></w:rPr>\n              <w:t>Git</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:rPr><w:lang w:val=\"en-US\"/></w:rPr>\n              <w:t xml:space=\"preserve\"> commit \\u2013amend (takes</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.start);\nawait context.sync();\n\n// The document only ever carries a single \"_GoBack\" bookmark; move it\n// from the end of the \"Files\" heading to the end of the paragraph we\n// just filled in.\ncontext.document.deleteBookmark(\"_GoBack\");\nconst endOfTarget = target.getRange(\"End\");\nendOfTarget.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Push\" heading paragraph (exact word, case-sensitive so we\n# skip the lower-case \"push\" occurrences inside other sentences), then\n# step back one paragraph to the empty list-item paragraph that sits\n# right above it \u2014 that's the paragraph the new text belongs in.\n$findRng = $d.Content\n$findRng.Find.Text = \"Push\"\n$findRng.Find.MatchCase = $true\n$findRng.Find.MatchWholeWord = $true\n$findRng.Find.Forward = $true\n$findRng.Find.Execute() | Out-Null\n\n$pushParagraph = $findRng.Paragraphs(1)\n$target = $pushParagraph.Previous()\n\n# Pull the paragraph's own opening <w:p ...> attributes and its\n# <w:pPr> block so the replacement keeps the exact same paragraph\n# (list/style/rPr) properties it already has.\n$xml = $target.Range.WordOpenXML\n\n$attrs = \"\"\nif ($xml -match '<w:p\\s+([^>]*)>') {\n  $attrs = $matches[1]\n}\n$attrs = [regex]::Replace($attrs, 'w14:\\w+=\"[^\"]*\"\\s*', '')\n\n$pPr = \"\"\nif ($xml -match '(?s)<w:pPr>.*?</w:pPr>') {\n  $pPr = $matches[0]\n}\n\n# Remove the document's single \"_GoBack\" bookmark from its current spot\n# (end of the \"Files\" heading) before we re-author it further down \u2014\n# this avoids a duplicate-id collision with the copy we are about to\n# insert.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$dash = [char]0x2013\n$newParagraphXml = '<w:p ' + $attrs + '>' + $pPr + `\n  '<w:proofErr w:type=\"spellStart\"/>' + `\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Git</w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellEnd\"/>' + `\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> commit ' + $dash + 'amend (takes</w:t></w:r>' + `\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' + `\n  '</w:p>'\n\n$ooxmlPackage = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $newParagraphXml + '</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>'\n\n$insertionPoint = $target.Range.Duplicate\n$insertionPoint.Collapse(1)\n$insertionPoint.InsertXML($ooxmlPackage)\n"}
